# Update workbook: Avverkningsanmälningar - refresh "Förändrad" (changed) date
# for all existing data rows, mark row 458 with explicit row height, and
# append the new entry "A 44611-2023" as row 459.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C ("Förändrad") from 45190 to 45192 for all data rows (2-458).
$ws.Range("C2:C458").Value = 45192

# 2. Row 458 now gets an explicit row height (matches the other data rows).
$ws.Rows.Item(458).RowHeight = 15

# 3. Append the new row 459 with the new case data.
$ws.Cells.Item(459,1).Value = "A 44611-2023"

$ws.Cells.Item(459,2).Value = 45189
$ws.Cells.Item(459,2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(459,3).Value = 45192
$ws.Cells.Item(459,3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(459,4).Value = "GÄVLEBORGS LÄN"
$ws.Cells.Item(459,5).Value = "OVANÅKER"

$ws.Cells.Item(459,7).Value = 1
$ws.Cells.Item(459,8).Value = 0
$ws.Cells.Item(459,9).Value = 0
$ws.Cells.Item(459,10).Value = 0
$ws.Cells.Item(459,11).Value = 0
$ws.Cells.Item(459,12).Value = 0
$ws.Cells.Item(459,13).Value = 0
$ws.Cells.Item(459,14).Value = 0
$ws.Cells.Item(459,15).Value = 0
$ws.Cells.Item(459,16).Value = 0
$ws.Cells.Item(459,17).Value = 0

$ws.Cells.Item(459,18).WrapText = $true
